$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 12.56082016290097
$ws.Range("C8").Value = 12.04303824191804
$ws.Range("D8").Value = 96.21070741913216
$ws.Range("E8").Value = 92.51766565866501
